$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.1298123333333333
$ws.Range("H2").Value = 0.389437
$ws.Range("I2").Value = 0.01442185502613333
$ws.Range("J2").Value = 0.01442185502613333
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 3.293530542941445
$ws.Range("R2").Value = 29.641774886473
$ws.Range("S2").Value = 0.008337103437859896
$ws.Range("T2").Value = 0.008337103437859898
$ws.Range("G3").Value = 0.1298123333333333
$ws.Range("H3").Value = 0.389437
$ws.Range("I3").Value = 0.01442185502613333
$ws.Range("J3").Value = 0.01442185502613333
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 1.326642021197222
$ws.Range("R3").Value = 11.939778190775
$ws.Range("S3").Value = 0.003358205309325837
$ws.Range("T3").Value = 0.003358205309325837
$ws.Range("G4").Value = 0.1298123333333333
$ws.Range("H4").Value = 0.389437
$ws.Range("I4").Value = 0.01442185502613333
$ws.Range("J4").Value = 0.01442185502613333
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 1.077108316262222
$ws.Range("R4").Value = 9.693974846360002
$ws.Range("S4").Value = 0.002726546278947597
$ws.Range("T4").Value = 0.002726546278947599
$ws.Range("I5").Value = 0.7859600471098795
$ws.Range("J5").Value = 0.7859600471098797
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 179.4903232626736
$ws.Range("R5").Value = 1615.412909364062
$ws.Range("S5").Value = 0.454354186677547
$ws.Range("T5").Value = 0.4543541866775471
$ws.Range("I6").Value = 0.7859600471098795
$ws.Range("J6").Value = 0.7859600471098797
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.1830149587788528
$ws.Range("T6").Value = 0.1830149587788528
$ws.Range("I7").Value = 0.7859600471098795
$ws.Range("J7").Value = 0.7859600471098797
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.1485909016534798
$ws.Range("T7").Value = 0.1485909016534798
$ws.Range("I8").Value = 0.1996180978639869
$ws.Range("J8").Value = 0.199618097863987
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 45.58694433188911
$ws.Range("R8").Value = 410.282498987002
$ws.Range("S8").Value = 0.1153968561565204
$ws.Range("T8").Value = 0.1153968561565204
$ws.Range("I9").Value = 0.1996180978639869
$ws.Range("J9").Value = 0.199618097863987
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("S9").Value = 0.04648213110377496
$ws.Range("T9").Value = 0.04648213110377497
$ws.Range("I10").Value = 0.1996180978639869
$ws.Range("J10").Value = 0.199618097863987
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("S10").Value = 0.03773911060369157
$ws.Range("T10").Value = 0.03773911060369158
